$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1)
$ws.Range("A1").Value = "Nama"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Username"
$ws.Range("D1").Value = "Password"
$ws.Range("E1").Value = "Peran"
$ws.Range("F1").Value = "Nomor Telepon"
$ws.Range("G1").Value = "Notifikasi WA"

# Row 2 - Aminu
$ws.Range("A2").Value = "aminu bil huda"
$ws.Range("B2").Value = "aminu@aminu.com"
$ws.Range("C2").Value = "aminu"
$ws.Range("D2").Value = 12345678
$ws.Range("E2").Value = "Guru"
$ws.Range("F2").Value = 85707357080
$ws.Range("G2").Value = "Ya"

# Row 3 - Danang
$ws.Range("A3").Value = "danang putra"
$ws.Range("B3").Value = "danang@danang.com"
$ws.Range("C3").Value = "danang"
$ws.Range("D3").Value = 12345678
$ws.Range("E3").Value = "Staf TU"
$ws.Range("F3").Value = 897648391
$ws.Range("G3").Value = "Tidak"

# Replace existing hyperlinks on email column with new addresses
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:aminu@aminu.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:danang@danang.com")
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"

# Column width for F (Nomor Telepon) - best-fit width ends up stored as 12
$ws.Columns.Item(6).ColumnWidth = 11.14

$ws.Range("E8").Select()
